$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- CoAP row (row 6): results were previously blank, now populated ---
$ws.Range("B6").Value = 597.81113861708502
$ws.Range("C6").Value = 2029.4212925735801
$ws.Range("D6").Value = 6318.0061176205299
$ws.Range("E6").Value = 26075.643014648202
$ws.Range("F6").Value = 5772.0237275871305
$ws.Range("G6").Value = 158178.70910074501
$ws.Range("H6").Value = 3513.3936214032301
$ws.Range("I6").Value = 10509.4771512315
$ws.Range("J6").Value = 0.38400000000000001
$ws.Range("K6").Value = 0.38400000000000001
$ws.Range("L6").Value = 0.38400000000000001
$ws.Range("M6").Value = 0.38400000000000001

# Entering the CoAP "Avg throughput" values across J6:M6 as one uniform
# block also drops the inner left borders so the row reads as a single
# continuous strip (matching J6's existing formatting).
$xlEdgeLeft = 7
$xlEdgeRight = 10
$xlThin = 2
$xlContinuous = 1
$xlLineStyleNone = -4142

$ws.Range("M6").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range("M6").Borders.Item($xlEdgeRight).Weight = $xlThin
$ws.Range("K6").Borders.Item($xlEdgeLeft).LineStyle = $xlLineStyleNone
$ws.Range("L6").Borders.Item($xlEdgeLeft).LineStyle = $xlLineStyleNone
$ws.Range("M6").Borders.Item($xlEdgeLeft).LineStyle = $xlLineStyleNone

# --- HTTP row (row 7): re-run results replace the earlier measurements ---
$ws.Range("B7").Value = 328.84254832621701
$ws.Range("C7").Value = 743.64511756416596
$ws.Range("D7").Value = 35722.251096703898
$ws.Range("E7").Value = 1038535.99451121
$ws.Range("F7").Value = 1876659.7296289001
$ws.Range("G7").Value = 45562676.805269897
$ws.Range("H7").Value = 3489190.86420677
$ws.Range("I7").Value = 21641500.804111298
$ws.Range("J7").Value = 0.8
$ws.Range("K7").Value = 81.92
$ws.Range("L7").Value = 8388.6080000000002
$ws.Range("M7").Value = 82561.296000000002

# Leave selection where the user last clicked after entering the data.
$ws.Range("O8").Select()
